$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "USE DataBaseName;"
$ws.Range("B5").Value = "switch to another DataBase"
$ws.Range("A6").Value = "SELECT NEWID();"
$ws.Range("B6").Value = "returns a guid (globally unique identifier). Return-type is uniqueidentifier"

$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.Range("B7").Select() | Out-Null
